$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("summary")

# --- Update existing "max" rows (81-86): relabel from "velocity..." to "max..." ---
# (order matters for shared-string table append order)
$ws.Range("A81").Value = "max, daily 2000, bio=100, dt=1 hr, deep"
$ws.Range("A82").Value = "max, daily 2000, bio=100, dt=1 hr, shallow"
$ws.Range("A83").Value = "max, daily 2000, bio=100, dt=1 hr, temp const"
$ws.Range("A84").Value = "max, daily 2000, bio=100, dt=1 hr, temp daily"
$ws.Range("A85").Value = "max, daily 2000, bio=100, dt=1 hr, zoop const"
$ws.Range("A86").Value = "max, daily 2000, bio=100, dt=1 hr, zoop daily"

# Row 83 gains a B value (mass-conservation residual) and an (empty, formatted) G cell
$ws.Range("B83").Value = [double]"-4.4668999999999999E-13"
$ws.Range("G83").NumberFormat = "0.00E+00"

# --- New "grad" rows (87-92) ---
$ws.Range("A87").Value = "grad, daily 2000, bio=100, dt=1 hr, deep"
$ws.Range("F87").Value = [double]"3.0599999999999999E-2"

$ws.Range("A89").Value = "grad, daily 2000, bio=100, dt=1 hr, temp const"
$ws.Range("A88").Value = "grad, daily 2000, bio=100, dt=1 hr, shallow"
$ws.Range("A90").Value = "grad, daily 2000, bio=100, dt=1 hr, temp daily"
$ws.Range("A91").Value = "grad, daily 2000, bio=100, dt=1 hr, zoop const"
$ws.Range("A92").Value = "grad, daily 2000, bio=100, dt=1 hr,zoop daily"

$ws.Range("F89").Value = [double]"-2.1473E-4"
$ws.Range("F89").NumberFormat = "0.00E+00"

$ws.Range("G89").Value = "Not enough change in distribution"
$ws.Range("G87").Value = "Unstable"

# --- View state: scroll position and active selection ---
$ws.Select()
$excel.ActiveWindow.ScrollRow = 60
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G87").Select()
